# Apply scraped crypto price/volume updates from coinranking.com
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'28.509.10"
$ws.Range('E2').Value = "  +0.11%  "

$ws.Range('D3').Value = "'1.827.35"
$ws.Range('E3').Value = "  -0.09%  "

$ws.Range('E4').Value = "  +0.28%  "

$ws.Range('D5').Value = "'316.46"
$ws.Range('E5').Value = "  +0.30%  "

$ws.Range('E6').Value = "  +0.24%  "

$ws.Range('D7').Value = "'0.5158"
$ws.Range('E7').Value = "  +1.96%  "

$ws.Range('D8').Value = "'0.3871"

$ws.Range('D9').Value = "'0.08315"
$ws.Range('E9').Value = "  +7.75%  "

$ws.Range('B10').Value = "Polygon"
$ws.Range('C10').Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range('D10').Value = "'1.120"
$ws.Range('E10').Value = "  +0.56%  "

$ws.Range('B11').Value = "OKB"
$ws.Range('C11').Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range('D11').Value = "'42.03"
$ws.Range('E11').Value = "  +0.31%  "

$ws.Range('D12').Value = "'6.417"
$ws.Range('E12').Value = "  +2.57%  "

$ws.Range('D13').Value = "'21.19"
$ws.Range('E13').Value = "  +0.73%  "

$ws.Range('D14').Value = "'1.004"
$ws.Range('E14').Value = "  +0.24%  "

$ws.Range('D15').Value = "'7.501"
$ws.Range('E15').Value = "  -0.94%  "

$ws.Range('D16').Value = "'1.823.28"
$ws.Range('E16').Value = "  -0.05%  "

$ws.Range('D17').Value = "'93.94"
$ws.Range('E17').Value = "  +0.53%  "

$ws.Range('D18').Value = "'0.00001124"
$ws.Range('E18').Value = "  +3.91%  "

$ws.Range('D19').Value = "'0.06650"
$ws.Range('E19').Value = "  +0.40%  "

$ws.Range('D20').Value = "'17.80"
$ws.Range('E20').Value = "  +0.48%  "

$ws.Range('E21').Value = "  +0.25%  "

$ws.Range('D22').Value = "'6.059"

$ws.Range('D23').Value = "'28.551.64"
$ws.Range('E23').Value = "  +0.15%  "

$ws.Range('D24').Value = "'11.43"
$ws.Range('E24').Value = "  +2.42%  "

$ws.Range('D25').Value = "'2.285"
$ws.Range('E25').Value = "  +1.05%  "

$ws.Range('D26').Value = "'21.18"
$ws.Range('E26').Value = "  +2.90%  "

$ws.Range('D27').Value = "'159.72"

$ws.Range('D28').Value = "'2.034.04"
$ws.Range('E28').Value = "  -0.03%  "

$ws.Range('D29').Value = "'2.406"
$ws.Range('E29').Value = "  -0.34%  "

$ws.Range('D30').Value = "'126.10"
$ws.Range('E30').Value = "  +0.73%  "

$ws.Range('D31').Value = "'0.1094"
$ws.Range('E31').Value = "  +0.71%  "

$ws.Range('D32').Value = "'1.096"
$ws.Range('E32').Value = "  -3.20%  "

$ws.Range('D33').Value = "'0.07596"

$ws.Range('D34').Value = "'5.739"
$ws.Range('E34').Value = "  +1.44%  "

$ws.Range('D35').Value = "'3.674"
$ws.Range('E35').Value = "  +0.25%  "

$ws.Range('D36').Value = "'0.2230"
$ws.Range('E36').Value = "  +0.43%  "

$ws.Range('D37').Value = "'0.02382"
$ws.Range('E37').Value = "  +2.47%  "

$ws.Range('D38').Value = "'5.262"
$ws.Range('E38').Value = "  +2.42%  "

$ws.Range('D39').Value = "'11.90"
$ws.Range('E39').Value = "  +6.24%  "

$ws.Range('D40').Value = "'8.761"
$ws.Range('E40').Value = "  -2.10%  "

$ws.Range('D41').Value = "'0.6379"
$ws.Range('E41').Value = "  +2.27%  "

$ws.Range('E42').Value = "  +0.26%  "

$ws.Range('D43').Value = "'1.396"
$ws.Range('E43').Value = "  -0.30%  "

$ws.Range('D44').Value = "'13.54"
$ws.Range('E44').Value = "  +0.30%  "

$ws.Range('D45').Value = "'0.6100"
$ws.Range('E45').Value = "  +3.46%  "

$ws.Range('D46').Value = "'3.804"
$ws.Range('E46').Value = "  +2.33%  "

$ws.Range('D47').Value = "'127.50"
$ws.Range('E47').Value = "  +2.37%  "

$ws.Range('D48').Value = "'1.996"

$ws.Range('E49').Value = "  +1.59%  "

$ws.Range('D50').Value = "'0.06982"
$ws.Range('E50').Value = "  +0.75%  "

$ws.Range('D51').Value = "'74.30"
$ws.Range('E51').Value = "  +0.63%  "
